$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "42.196.99"
$ws.Range("E2").Value = "  -1.67%  "
$ws.Range("D3").Value = "2.295.51"
$ws.Range("E3").Value = "  -2.97%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue $ws.Range("D5") "317.54"
$ws.Range("E5").Value = "  -1.13%  "
Set-TextValue $ws.Range("D6") "103.80"
$ws.Range("E6").Value = "  -3.32%  "
$ws.Range("E7").Value = "  -1.01%  "
$ws.Range("E8").Value = "  +0.13%  "
Set-TextValue $ws.Range("D9") "0.606"
$ws.Range("E9").Value = "  -2.61%  "
Set-TextValue $ws.Range("D10") "39.50"
$ws.Range("E10").Value = "  -5.51%  "
Set-TextValue $ws.Range("D11") "0.0906"
$ws.Range("E11").Value = "  -2.44%  "
Set-TextValue $ws.Range("D12") "8.28"
$ws.Range("E12").Value = "  -3.12%  "
$ws.Range("E13").Value = "  -0.09%  "
Set-TextValue $ws.Range("D14") "0.963"
$ws.Range("E14").Value = "  -4.54%  "
Set-TextValue $ws.Range("D15") "15.27"
$ws.Range("E15").Value = "  -5.39%  "
$ws.Range("D16").Value = "2.643.24"
$ws.Range("E16").Value = "  -3.05%  "
$ws.Range("D17").Value = "2.298.46"
$ws.Range("E17").Value = "  -5.37%  "
$ws.Range("D18").Value = "42.161.82"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("E19").Value = "  -2.08%  "
Set-TextValue $ws.Range("D20") "0.0000106"
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D21") "73.40"
$ws.Range("E21").Value = "  -3.78%  "
$ws.Range("B22").Value = "PancakeSwap"
$ws.Range("C22").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D22") "3.60"
$ws.Range("E22").Value = "  -3.14%  "
Set-TextValue $ws.Range("D23") "278.63"
$ws.Range("E23").Value = "  +5.49%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D24") "2.27"
$ws.Range("E24").Value = "  -2.74%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D25") "9.95"
$ws.Range("E25").Value = "  +5.10%  "
$ws.Range("E26").Value = "  +0.55%  "
Set-TextValue $ws.Range("D27") "10.84"
$ws.Range("E27").Value = "  -5.46%  "
Set-TextValue $ws.Range("D28") "2.33"
$ws.Range("E28").Value = "  +3.36%  "
Set-TextValue $ws.Range("D29") "22.93"
$ws.Range("E29").Value = "  -1.35%  "
Set-TextValue $ws.Range("D30") "36.68"
$ws.Range("E30").Value = "  -0.92%  "
Set-TextValue $ws.Range("D31") "163.17"
$ws.Range("E31").Value = "  -4.75%  "
Set-TextValue $ws.Range("D32") "0.0874"
$ws.Range("E32").Value = "  -2.99%  "
Set-TextValue $ws.Range("D33") "2.85"
$ws.Range("E33").Value = "  -3.11%  "
Set-TextValue $ws.Range("D34") "5.83"
$ws.Range("E34").Value = "  -3.17%  "
Set-TextValue $ws.Range("D35") "0.136"
$ws.Range("E35").Value = "  +3.59%  "
$ws.Range("E36").Value = "  -7.00%  "
Set-TextValue $ws.Range("D37") "4.55"
$ws.Range("E37").Value = "  -3.86%  "
Set-TextValue $ws.Range("D38") "2.94"
$ws.Range("E38").Value = "  +8.38%  "
Set-TextValue $ws.Range("D39") "0.0351"
$ws.Range("E39").Value = "  -3.89%  "
$ws.Range("E40").Value = "  -4.62%  "
Set-TextValue $ws.Range("D41") "99.05"
$ws.Range("E41").Value = "  -0.33%  "
Set-TextValue $ws.Range("D42") "1.45"
$ws.Range("E42").Value = "  -4.97%  "
Set-TextValue $ws.Range("D43") "69.48"
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("E44").Value = "  +0.14%  "
Set-TextValue $ws.Range("D45") "0.225"
$ws.Range("E45").Value = "  -6.68%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D46") "113.58"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D47") "11.96"
$ws.Range("E47").Value = "  -2.74%  "
Set-TextValue $ws.Range("D48") "77.13"
$ws.Range("E48").Value = "  -0.25%  "
Set-TextValue $ws.Range("D49") "9.00"
$ws.Range("E49").Value = "  -2.61%  "
Set-TextValue $ws.Range("D50") "5.29"
$ws.Range("E50").Value = "  -4.89%  "
$ws.Range("D51").Value = "1.582.49"
$ws.Range("E51").Value = "  +0.42%  "
